$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.2475347124312
$ws.Range("C2").Value = 8.263854200641449
$ws.Range("D2").Value = 3.886335313579849
$ws.Range("F2").Value = 24.07170157423461
$ws.Range("G2").Value = 29.52698303746906
$ws.Range("H2").Value = 13.65507026307969
$ws.Range("M2").Value = 19.96853898132113

$ws.Range("B3").Value = 10.71492775045894
$ws.Range("C3").Value = 7.922581529107058
$ws.Range("D3").Value = 3.892723106426581
$ws.Range("F3").Value = 23.82015167278259
$ws.Range("G3").Value = 28.96322266574675
$ws.Range("H3").Value = 13.64876489467059
$ws.Range("M3").Value = 19.36414585261561

$ws.Range("B4").Value = 10.37589454253218
$ws.Range("C4").Value = 7.703418244744816
$ws.Range("D4").Value = 3.896783410241674
$ws.Range("F4").Value = 23.67348880842128
$ws.Range("G4").Value = 28.62493826143402
$ws.Range("H4").Value = 13.64855803436505
$ws.Range("M4").Value = 18.99168008987713

$ws.Range("B5").Value = 10.23491323021945
$ws.Range("C5").Value = 7.611759913204648
$ws.Range("D5").Value = 3.898473348474039
$ws.Range("F5").Value = 23.61574415374162
$ws.Range("G5").Value = 28.48927497503779
$ws.Range("H5").Value = 13.6493928279958
$ws.Range("M5").Value = 18.83980783706242

$ws.Range("B6").Value = 10.21133927762106
$ws.Range("C6").Value = 7.596400634200308
$ws.Range("D6").Value = 3.898756112672826
$ws.Range("F6").Value = 23.60627953415904
$ws.Range("G6").Value = 28.46688659580859
$ws.Range("H6").Value = 13.64958687157188
$ws.Range("M6").Value = 18.81459182937592

$ws.Range("B7").Value = 10.37400436462951
$ws.Range("C7").Value = 7.7021915110846
$ws.Range("D7").Value = 3.896806057500229
$ws.Range("F7").Value = 23.67270177899914
$ws.Range("G7").Value = 28.62309951298588
$ws.Range("H7").Value = 13.64856557499311
$ws.Range("M7").Value = 18.98963190478101

$ws.Range("B8").Value = 11.06649417296349
$ws.Range("C8").Value = 8.148219911062412
$ws.Range("D8").Value = 3.888509497530909
$ws.Range("F8").Value = 23.98338967033638
$ws.Range("G8").Value = 29.3310941698613
$ws.Range("H8").Value = 13.65213431839109
$ws.Range("M8").Value = 19.76057610694909

$ws.Range("B9").Value = 12.32189773246078
$ws.Range("C9").Value = 8.943902716147695
$ws.Range("D9").Value = 3.873310264554926
$ws.Range("F9").Value = 24.65145310088139
$ws.Range("G9").Value = 30.77242873923063
$ws.Range("H9").Value = 13.68828871363744
$ws.Range("M9").Value = 21.25110361432481

$ws.Range("B10").Value = 13.17370492927803
$ws.Range("C10").Value = 9.477424420112531
$ws.Range("D10").Value = 3.862760399131219
$ws.Range("F10").Value = 25.17381652207041
$ws.Range("G10").Value = 31.8510136040202
$ws.Range("H10").Value = 13.73269027763679
$ws.Range("M10").Value = 22.32034253671616

$ws.Range("B11").Value = 13.5446277780731
$ws.Range("C11").Value = 9.708583484889111
$ws.Range("D11").Value = 3.858087510881728
$ws.Range("F11").Value = 25.41731968459154
$ws.Range("G11").Value = 32.34334116611193
$ws.Range("H11").Value = 13.7567619499317
$ws.Range("M11").Value = 22.79877583828171

$ws.Range("B12").Value = 13.68262117814696
$ws.Range("C12").Value = 9.794426884345887
$ws.Range("D12").Value = 3.856335583375187
$ws.Range("F12").Value = 25.51028813666762
$ws.Range("G12").Value = 32.5298046883828
$ws.Range("H12").Value = 13.76643280153437
$ws.Range("M12").Value = 22.97862599430049

$ws.Range("B13").Value = 13.65301274425614
$ws.Range("C13").Value = 9.776014651329518
$ws.Range("D13").Value = 3.856712118781648
$ws.Range("F13").Value = 25.49023333632165
$ws.Range("G13").Value = 32.4896483735222
$ws.Range("H13").Value = 13.76432533992179
$ws.Range("M13").Value = 22.93995346730165

$ws.Range("B14").Value = 13.55603045357605
$ws.Range("C14").Value = 9.715679936207716
$ws.Range("D14").Value = 3.857943029376946
$ws.Range("F14").Value = 25.42495353809243
$ws.Range("G14").Value = 32.35868221280798
$ws.Range("H14").Value = 13.75754645645891
$ws.Range("M14").Value = 22.813599663564

$ws.Range("B15").Value = 13.49630246102696
$ws.Range("C15").Value = 9.678502045940055
$ws.Range("D15").Value = 3.85869927214525
$ws.Range("F15").Value = 25.38506406186352
$ws.Range("G15").Value = 32.2784594223692
$ws.Range("H15").Value = 13.753466470035
$ws.Range("M15").Value = 22.73602719781077

$ws.Range("B16").Value = 13.1491232852912
$ws.Range("C16").Value = 9.462082529151756
$ws.Range("D16").Value = 3.863068276144054
$ws.Range("F16").Value = 25.158013952211
$ws.Range("G16").Value = 31.81885429509699
$ws.Range("H16").Value = 13.73119496220733
$ws.Range("M16").Value = 22.28889997214678

$ws.Range("B17").Value = 12.93183159981898
$ws.Range("C17").Value = 9.326335500926517
$ws.Range("D17").Value = 3.865780460727885
$ws.Range("F17").Value = 25.02017052766331
$ws.Range("G17").Value = 31.5371774025777
$ws.Range("H17").Value = 13.71852328076846
$ws.Range("M17").Value = 22.0124298647587

$ws.Range("B18").Value = 12.80529568890266
$ws.Range("C18").Value = 9.247172473937669
$ws.Range("D18").Value = 3.86735235892036
$ws.Range("F18").Value = 24.94144509063851
$ws.Range("G18").Value = 31.37533398585553
$ws.Range("H18").Value = 13.71159956624242
$ws.Range("M18").Value = 21.85267090783623

$ws.Range("B19").Value = 12.7621884682957
$ws.Range("C19").Value = 9.220183952303023
$ws.Range("D19").Value = 3.867886641525253
$ws.Range("F19").Value = 24.91488853033221
$ws.Range("G19").Value = 31.32057229542021
$ws.Range("H19").Value = 13.70931798739232
$ws.Range("M19").Value = 21.7984577207549

$ws.Range("B20").Value = 12.95512431591213
$ws.Range("C20").Value = 9.340898536506204
$ws.Range("D20").Value = 3.86549051436763
$ws.Range("F20").Value = 25.03478701789123
$ws.Range("G20").Value = 31.56714638276706
$ws.Range("H20").Value = 13.71983446634611
$ws.Range("M20").Value = 22.0419385683425

$ws.Range("B21").Value = 13.58458405925349
$ws.Range("C21").Value = 9.733447838099021
$ws.Range("D21").Value = 3.857581008097329
$ws.Range("F21").Value = 25.44410788626506
$ws.Range("G21").Value = 32.39715097861318
$ws.Range("H21").Value = 13.75952252048997
$ws.Range("M21").Value = 22.85075003856932

$ws.Range("B22").Value = 13.98156326325118
$ws.Range("C22").Value = 9.980130949339236
$ws.Range("D22").Value = 3.85251397965713
$ws.Range("F22").Value = 25.71600687388121
$ws.Range("G22").Value = 32.93968906354509
$ws.Range("H22").Value = 13.78869720156525
$ws.Range("M22").Value = 23.37157916598345

$ws.Range("B23").Value = 13.77103068253317
$ws.Range("C23").Value = 9.849383806978922
$ws.Range("D23").Value = 3.855209174471347
$ws.Range("F23").Value = 25.57051675553482
$ws.Range("G23").Value = 32.65018611786212
$ws.Range("H23").Value = 13.77283070705358
$ws.Range("M23").Value = 23.09436804819591

$ws.Range("B24").Value = 12.94459869350832
$ws.Range("C24").Value = 9.334318074994748
$ws.Range("D24").Value = 3.865621559743161
$ws.Range("F24").Value = 25.02817726840984
$ws.Range("G24").Value = 31.55359709051584
$ws.Range("H24").Value = 13.71924055354745
$ws.Range("M24").Value = 22.02860020098348

$ws.Range("B25").Value = 11.99417181612458
$ws.Range("C25").Value = 8.737436798342815
$ws.Range("D25").Value = 3.877311170978883
$ws.Range("F25").Value = 24.46486037857789
$ws.Range("G25").Value = 30.37815568248895
$ws.Range("H25").Value = 13.67537705923278
$ws.Range("M25").Value = 20.85155013784399
